$wb = $excel.ActiveWorkbook

# --- Update conversion message on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.36 = 49753.03 pesos`n✅ 49753.03 pesos = 12.35 = 967.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 80.90000000000001
$wsTasas.Range("O10").Value = 4025.02
$wsTasas.Range("N12").Value = 4029.99
$wsTasas.Range("O12").Value = 78.40600000000001
